{"js": "// Add a new \"Problem 7\" row to the Solutions table, matching the\n// existing table's pattern (Problem | Part | Solution):\n//   7 | - | 0.377\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[tables.items.length - 1];\n\n// Append one row at the end of the table, populating all three\n// columns in a single call so the write lands on the newly created\n// row (not the existing header/last row).\ntable.addRows(\"End\", 1, [[\"7\", \"-\", \"0.377\"]]);\n\nawait context.sync();\n", "ps1": "# Add a new \"Problem 7\" row to the Solutions table, matching the\n# existing table's pattern (Problem | Part | Solution):\n#   7 | - | 0.377\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item($d.Tables.Count)\n\n$newRow = $t.Rows.Add()\n$idx = $newRow.Index\n\n$t.Cell($idx, 1).Range.Text = \"7\"\n$t.Cell($idx, 2).Range.Text = \"-\"\n$t.Cell($idx, 3).Range.Text = \"0.377\"\n"}
